$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update row 9 (Ano 2025) with refreshed faturamento figures
$ws.Range("B9").Value = 4037834.59
$ws.Range("C9").Value = 634141.24
$ws.Range("D9").Value = 4671975.83
$ws.Range("E9").Value = 13.57329881563193
$ws.Range("F9").Value = 86.42670118436807
$ws.Range("G9").Value = -38.71341112799489
$ws.Range("H9").Value = -27.08228855952724
$ws.Range("I9").Value = 40587
$ws.Range("J9").Value = 1742
$ws.Range("K9").Value = 42329
$ws.Range("L9").Value = 29304
$ws.Range("M9").Value = 159.4313346300846
$ws.Range("N9").Value = 8.847019406241952
